# Applies the cryptos-list price/volume refresh described by the commit diff.
# Column D holds price text, column E holds percentage-change text; both are
# stored as plain strings (inline strings) in the workbook, not numbers.
# For D-column values that look like plain numbers (e.g. '7.80', '18.50'),
# Excel's Range.Value setter would auto-convert them to numeric, dropping
# trailing zeros / changing type. Prefixing with a leading apostrophe forces
# Excel to keep them as literal text, exactly like typing '7.80 into a cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.797.80'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").Value = '2.907.38'
$ws.Range("E3").Value = '  -2.79%  '
$ws.Range("D5").Value = '''525.99'
$ws.Range("E5").Value = '  -2.97%  '
$ws.Range("D6").Value = '''144.14'
$ws.Range("E6").Value = '  -4.92%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -4.31%  '
$ws.Range("D9").Value = '2.915.45'
$ws.Range("E10").Value = '  -5.34%  '
$ws.Range("D11").Value = '''6.16'
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("E12").Value = '  -2.97%  '
$ws.Range("D13").Value = '3.415.80'
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").Value = '60.820.00'
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("E16").Value = '  -6.05%  '
$ws.Range("D17").Value = '2.923.53'
$ws.Range("E17").Value = '  -2.49%  '
$ws.Range("E18").Value = '  -3.96%  '
$ws.Range("E19").Value = '  -5.20%  '
$ws.Range("E20").Value = '  -4.29%  '
$ws.Range("D21").Value = '''353.31'
$ws.Range("E21").Value = '  -6.49%  '
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '''64.84'
$ws.Range("E25").Value = '  -2.02%  '
$ws.Range("E26").Value = '  -4.09%  '
$ws.Range("E27").Value = '  -4.94%  '
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '0.0₃0866'
$ws.Range("E29").Value = '  -7.29%  '
$ws.Range("D30").Value = '''7.80'
$ws.Range("E30").Value = '  -5.37%  '
$ws.Range("D31").Value = '''0.999'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '''1.67'
$ws.Range("E32").Value = '  -3.19%  '
$ws.Range("D33").Value = '''19.65'
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("D34").Value = '''153.37'
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("D35").Value = '''4.39'
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("D36").Value = '''5.58'
$ws.Range("E36").Value = '  -6.60%  '
$ws.Range("D37").Value = '''0.994'
$ws.Range("E37").Value = '  -6.96%  '
$ws.Range("E38").Value = '  -6.32%  '
$ws.Range("D39").Value = '''37.58'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  -5.44%  '
$ws.Range("E41").Value = '  -2.97%  '
$ws.Range("D42").Value = '2.291.30'
$ws.Range("E42").Value = '  -5.15%  '
$ws.Range("D43").Value = '''3.69'
$ws.Range("E43").Value = '  -5.29%  '
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("D45").Value = '''20.36'
$ws.Range("E45").Value = '  -7.66%  '
$ws.Range("E47").Value = '  -4.89%  '
$ws.Range("E48").Value = '  -3.23%  '
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("E50").Value = '  -4.24%  '
$ws.Range("D51").Value = '''18.50'
$ws.Range("E51").Value = '  -5.77%  '
